$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 110.3882446666667
$ws.Range("H2").Value = 331.164734
$ws.Range("I2").Value = 0.4177264991141899
$ws.Range("J2").Value = 0.4177264991141899
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3252056666666667
$ws.Range("N2").Value = 0.975617
$ws.Range("O2").Value = 0.0158278498560244
$ws.Range("P2").Value = 0.0158278498560244
$ws.Range("Q2").Value = 35.89888269898644
$ws.Range("R2").Value = 323.089944290878
$ws.Range("S2").Value = 0.006611712308862106
$ws.Range("T2").Value = 0.006611712308862106

# Row 3
$ws.Range("G3").Value = 110.3882446666667
$ws.Range("H3").Value = 331.164734
$ws.Range("I3").Value = 0.4177264991141899
$ws.Range("J3").Value = 0.4177264991141899
$ws.Range("O3").Value = 0.8133441666880411
$ws.Range("P3").Value = 0.8133441666880411
$ws.Range("Q3").Value = 1844.732360960921
$ws.Range("R3").Value = 16602.59124864829
$ws.Range("S3").Value = 0.3397554113255435
$ws.Range("T3").Value = 0.3397554113255435

# Row 4
$ws.Range("G4").Value = 110.3882446666667
$ws.Range("H4").Value = 331.164734
$ws.Range("I4").Value = 0.4177264991141899
$ws.Range("J4").Value = 0.4177264991141899
$ws.Range("M4").Value = 3.509903666666667
$ws.Range("N4").Value = 10.529711
$ws.Range("O4").Value = 0.1708279834559346
$ws.Range("P4").Value = 0.1708279834559346
$ws.Range("Q4").Value = 387.4521047124304
$ws.Range("R4").Value = 3487.068942411874
$ws.Range("S4").Value = 0.0713593754797843
$ws.Range("T4").Value = 0.0713593754797843

# Row 5
$ws.Range("I5").Value = 0.4902812436402899
$ws.Range("J5").Value = 0.4902812436402899
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3252056666666667
$ws.Range("N5").Value = 0.975617
$ws.Range("O5").Value = 0.0158278498560244
$ws.Range("P5").Value = 0.0158278498560244
$ws.Range("Q5").Value = 42.13414493042412
$ws.Range("R5").Value = 379.207304373817
$ws.Range("S5").Value = 0.007760097911563426
$ws.Range("T5").Value = 0.007760097911563425

# Row 6
$ws.Range("I6").Value = 0.4902812436402899
$ws.Range("J6").Value = 0.4902812436402899
$ws.Range("O6").Value = 0.8133441666880411
$ws.Range("P6").Value = 0.8133441666880411
$ws.Range("S6").Value = 0.3987673895513881
$ws.Range("T6").Value = 0.398767389551388

# Row 7
$ws.Range("I7").Value = 0.4902812436402899
$ws.Range("J7").Value = 0.4902812436402899
$ws.Range("M7").Value = 3.509903666666667
$ws.Range("N7").Value = 10.529711
$ws.Range("O7").Value = 0.1708279834559346
$ws.Range("P7").Value = 0.1708279834559346
$ws.Range("Q7").Value = 454.7485020755902
$ws.Range("R7").Value = 4092.736518680312
$ws.Range("S7").Value = 0.08375375617733848
$ws.Range("T7").Value = 0.08375375617733846

# Row 8
$ws.Range("G8").Value = 24.30983866666667
$ws.Range("H8").Value = 72.92951600000001
$ws.Range("I8").Value = 0.09199225724552029
$ws.Range("J8").Value = 0.09199225724552029
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.3252056666666667
$ws.Range("N8").Value = 0.975617
$ws.Range("O8").Value = 0.0158278498560244
$ws.Range("P8").Value = 0.0158278498560244
$ws.Range("Q8").Value = 7.905697290152445
$ws.Range("R8").Value = 71.151275611372
$ws.Range("S8").Value = 0.001456039635598868
$ws.Range("T8").Value = 0.001456039635598868

# Row 9
$ws.Range("G9").Value = 24.30983866666667
$ws.Range("H9").Value = 72.92951600000001
$ws.Range("I9").Value = 0.09199225724552029
$ws.Range("J9").Value = 0.09199225724552029
$ws.Range("O9").Value = 0.8133441666880411
$ws.Range("P9").Value = 0.8133441666880411
$ws.Range("Q9").Value = 406.2492905250512
$ws.Range("R9").Value = 3656.24361472546
$ws.Range("S9").Value = 0.07482136581110961
$ws.Range("T9").Value = 0.07482136581110961

# Row 10
$ws.Range("G10").Value = 24.30983866666667
$ws.Range("H10").Value = 72.92951600000001
$ws.Range("I10").Value = 0.09199225724552029
$ws.Range("J10").Value = 0.09199225724552029
$ws.Range("M10").Value = 3.509903666666667
$ws.Range("N10").Value = 10.529711
$ws.Range("O10").Value = 0.1708279834559346
$ws.Range("P10").Value = 0.1708279834559346
$ws.Range("Q10").Value = 85.32519187220845
$ws.Range("R10").Value = 767.9267268498761
$ws.Range("S10").Value = 0.01571485179881182
$ws.Range("T10").Value = 0.01571485179881182
